$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<kilo>"
$ws.Range("C2").Value = 24

# Row 3
$ws.Range("B3").Value = "<are>"
$ws.Range("C3").Value = 30

# Row 4
$ws.Range("C4").Value = 30

# Row 5
$ws.Range("B5").Value = "<long>"
$ws.Range("C5").Value = 41

# Row 6
$ws.Range("B6").Value = "<para>"
$ws.Range("C6").Value = 30

# Row 7
$ws.Range("C7").Value = 35

# Row 9
$ws.Range("B9").Value = "<from>"
$ws.Range("C9").Value = 30

# Row 10
$ws.Range("C10").Value = 33

# Row 11
$ws.Range("B11").Value = "<would>"
$ws.Range("C11").Value = 32

# Row 12
$ws.Range("C12").Value = 35

# Row 13
$ws.Range("B13").Value = "<nine>"
$ws.Range("C13").Value = 31

# Row 14
$ws.Range("C14").Value = 36

# Row 15
$ws.Range("B15").Value = "<cad>"
$ws.Range("C15").Value = 11
